$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column B values (row 2 through row 8)
$ws.Range("B2").Value = 92.60959625244141
$ws.Range("B3").Value = 70.54509735107422
$ws.Range("B4").Value = 23.50209999084473
$ws.Range("B5").Value = -4.403299808502197
$ws.Range("B6").Value = -110.4229965209961
$ws.Range("B7").Value = -81.25930023193359
$ws.Range("B8").Value = 70.74289703369141

# Add new row 9 with pasture data
$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 61.31409454345703
$ws.Range("C9").Value = 0.0083
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.375
$ws.Range("F9").Value = -0.2207999974489212
